$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values (stored as date serials, formatted with numFmtId 14)
$ws.Range("A2").Value = 46054
$ws.Range("B2").Value = 46060

# Update the active selection to match the new cursor position
$ws.Range("F11").Select()
